$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Summary": update capital / P&L / trade counters for trade #27
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.34   # Current Capital
$wsSummary.Range("B4").Value = 0.34      # Total P&L $
$wsSummary.Range("B6").Value = 27        # Total Trades
$wsSummary.Range("B7").Value = 12        # Winning Trades
$wsSummary.Range("B9").Value = 44.44     # Win Rate %

# ---------------------------------------------------------------------
# Sheet "Strategy Status": update MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.34
$wsStatus.Range("D4").Value = 27
$wsStatus.Range("E4").Value = 0.34
$wsStatus.Range("F4").Value = 0.34
$wsStatus.Range("G4").Value = 44.44

# ---------------------------------------------------------------------
# Append newly-closed trade #27 to both "All Trades" and "MarketMaking"
# sheets (row 28); these two sheets mirror the same trade log.
# ---------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(28, 1).Value = 27

    # Column B holds a plain text date string ("2026-02-17"); without
    # the quote-prefix Excel would auto-convert it into a date serial.
    $ws.Cells.Item(28, 2).Value = "'2026-02-17"

    $ws.Cells.Item(28, 3).Value = "04:08:57"
    $ws.Cells.Item(28, 4).Value = "MarketMaking"
    $ws.Cells.Item(28, 5).Value = "UP"
    $ws.Cells.Item(28, 6).Value = 0.5600000000000001
    $ws.Cells.Item(28, 7).Value = 0.58
    $ws.Cells.Item(28, 8).Value = "CLOSED"
    $ws.Cells.Item(28, 9).Value = 3.5714
    $ws.Cells.Item(28, 10).Value = 0.02
    $ws.Cells.Item(28, 11).Value = 100.34
    $ws.Cells.Item(28, 12).Value = 0
    $ws.Cells.Item(28, 13).Value = 0
    $ws.Cells.Item(28, 14).Value = 0.6
    $ws.Cells.Item(28, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(28, 16).Value = "early_exit"
    $ws.Cells.Item(28, 17).Value = 0.12
}
